$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 257, pushing existing rows 257-304
# down to 259-306.
$ws.Rows.Item(257).Insert()
$ws.Rows.Item(257).Insert()

# Populate new row 257.
$ws.Range("A257").Value = 4
$ws.Range("B257").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C257").Value = "Los Lagos"
$ws.Range("D257").Value = 44476
$ws.Range("E257").Value = 10
$ws.Range("F257").Value = 100112004
$ws.Range("G257").Value = "Cebolla"
$ws.Range("H257").Value = "Morada(o)"
$ws.Range("I257").Value = "1a (cosecha)"
$ws.Range("J257").Value = 120
$ws.Range("K257").Value = 12000
$ws.Range("L257").Value = 12000
$ws.Range("M257").Value = 12000
$ws.Range("N257").Value = "`$/malla 18 kilos"
$ws.Range("O257").Value = "Región de Arica y Parinacota"
$ws.Range("P257").Value = 667
$ws.Range("Q257").Value = 18
$ws.Range("R257").Value = "Hortaliza"

# Populate new row 258.
$ws.Range("A258").Value = 4
$ws.Range("B258").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C258").Value = "Los Lagos"
$ws.Range("D258").Value = 44476
$ws.Range("E258").Value = 10
$ws.Range("F258").Value = 100112004
$ws.Range("G258").Value = "Cebolla"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 300
$ws.Range("K258").Value = 9000
$ws.Range("L258").Value = 9000
$ws.Range("M258").Value = 9000
$ws.Range("N258").Value = "`$/malla 18 kilos"
$ws.Range("O258").Value = "Perú"
$ws.Range("P258").Value = 500
$ws.Range("Q258").Value = 18
$ws.Range("R258").Value = "Hortaliza"
